# Update the four sheets (신화인터텍, 드림시큐리티, 대원미디어, 성호전자)
# Each sheet has a daily date/amount table in columns A:B.
# Rows 100 and 101 (dates 45959 / 45960) had placeholder 0 values that are
# now filled with the real remn_amt, and a new row 102 (date 45961) is
# appended with a placeholder 0 value, same as the rows used to have.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "신화인터텍";   B100 = 536;  B101 = 529  },
    @{ Sheet = "드림시큐리티"; B100 = 5994; B101 = 5787 },
    @{ Sheet = "대원미디어";   B100 = 2697; B101 = 2612 },
    @{ Sheet = "성호전자";     B100 = 1409; B101 = 1404 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    # Fill in the previously-placeholder values for the last two existing rows.
    $ws.Range("B100").Value = $u.B100
    $ws.Range("B101").Value = $u.B101

    # Append the new row for the next date, with a placeholder 0 amount,
    # matching the style used by the existing date column cells.
    $ws.Range("A102").Value = 45961
    $ws.Range("A102").Style = $ws.Range("A101").Style
    $ws.Range("A102").NumberFormat = $ws.Range("A101").NumberFormat
    $ws.Range("B102").Value = 0
}
